$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 936
$ws.Range("K3").Value = 868
$ws.Range("J4").Value = 1794
$ws.Range("K4").Value = 200
$ws.Range("K5").Value = 52
$ws.Range("J6").Value = 11058
$ws.Range("K6").Value = 1220
$ws.Range("J7").Value = 29253
$ws.Range("K7").Value = 3276

# Norwood Park
$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 10

# Logan Square
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 39

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 52
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 195

# South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 60

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 134

# West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 14
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 56

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 32
$ws.Range("K7").Value = 100

# New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 21
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 91

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K4").Value = 13
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 97
$ws.Range("K8").Value = 195
$ws.Range("K10").Value = 20
$ws.Range("K14").Value = 18
$ws.Range("K18").Value = 27
$ws.Range("K19").Value = 89
$ws.Range("K20").Value = 79
$ws.Range("K21").Value = 10
$ws.Range("K27").Value = 40
$ws.Range("K32").Value = 6
$ws.Range("K33").Value = 134
$ws.Range("K34").Value = 22
$ws.Range("K35").Value = 6
$ws.Range("K36").Value = 38
$ws.Range("K37").Value = 100
$ws.Range("K42").Value = 103
$ws.Range("K43").Value = 32
$ws.Range("K44").Value = 28
$ws.Range("K46").Value = 6
$ws.Range("K47").Value = 24
$ws.Range("K48").Value = 32
$ws.Range("K50").Value = 21
$ws.Range("K53").Value = 39
$ws.Range("K54").Value = 58
$ws.Range("K59").Value = 9
$ws.Range("J63").Value = 88
$ws.Range("K63").Value = 10
$ws.Range("K64").Value = 19
$ws.Range("K65").Value = 91
$ws.Range("K67").Value = 140
$ws.Range("K69").Value = 10
$ws.Range("K70").Value = 8
$ws.Range("K73").Value = 39
$ws.Range("K76").Value = 47
$ws.Range("K77").Value = 21
$ws.Range("K78").Value = 45
$ws.Range("K79").Value = 92
$ws.Range("K82").Value = 4
$ws.Range("K83").Value = 60
$ws.Range("K85").Value = 165
$ws.Range("K89").Value = 52
$ws.Range("K91").Value = 34
$ws.Range("K92").Value = 14
$ws.Range("K95").Value = 56
$ws.Range("K96").Value = 53
$ws.Range("K97").Value = 25
$ws.Range("K98").Value = 20
$ws.Range("J101").Value = 29253
$ws.Range("K101").Value = 3276

# North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 41
$ws.Range("K3").Value = 39
$ws.Range("K4").Value = 9
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 140

# Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 58

# Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 32

# Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 28
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 89

# Irving Park
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 28

# River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 8
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 47

# Bridgeport
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 18

# Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 37
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 103

# Avondale
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 20

# Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 45

# Jefferson Park
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K6").Value = 1
$ws.Range("K7").Value = 6

# West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 19
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 53

# Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 34

# Chinatown
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 10

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 32
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 92

# Near South Side
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 19

# Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 17
$ws.Range("K3").Value = 24
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 79

# Calumet Heights
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 7
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 27

# Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 38

# Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 97

# Garfield Ridge
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 22

# Kenwood
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 10
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 24

# Wicker Park
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 20

# Lincoln Square
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 21

# Gold Coast
$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 6

# Portage Park
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 12
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 39

# Montclare
$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 9

# West Town
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 25

# West Elsdon
$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 14

# O'Hare
$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 8

# Galewood
$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 6

# Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 52

# Armour Square
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 6

# Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 40

# Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 32

# South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 65
$ws.Range("K3").Value = 49
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 165

# Sheffield & DePaul
$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K5").Value = 4
$ws.Range("K6").Value = 4

# Riverdale
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 21

# Archer Heights
$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 13
